$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in D23:D26 to 844
$ws.Range("D23").Value = 844
$ws.Range("D24").Value = 844
$ws.Range("D25").Value = 844
$ws.Range("D26").Value = 844
